$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 155
$ws1.Range("F6").Value = 5139
$ws1.Range("F7").Value = 111
$ws1.Range("F8").Value = 5297
$ws1.Range("F9").Value = 609
$ws1.Range("F10").Value = 1345
$ws1.Range("F11").Value = 101

# Sheet "全部类型" (index 4): update column F (想去人数) values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 155
$ws4.Range("F7").Value = 5139
$ws4.Range("F8").Value = 111
$ws4.Range("F9").Value = 5297
$ws4.Range("F10").Value = 609
$ws4.Range("F11").Value = 1345
$ws4.Range("F12").Value = 101
